# Insert a new weekly price record at row 395 of the "Berenjena" sheet.
# This shifts the existing rows 395:405 down to 396:406 (dimension grows
# from A1:R405 to A1:R406), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 395; rows below shift down
# and inherit the formatting (incl. the date number format in column D).
$ws.Rows.Item(395).Insert()

# Populate the newly inserted row 395 with the new weekly observation.
$ws.Cells.Item(395, 1).Value  = 9
$ws.Cells.Item(395, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(395, 3).Value  = "Metropolitana"
$ws.Cells.Item(395, 4).Value  = 45239
$ws.Cells.Item(395, 5).Value  = 13
$ws.Cells.Item(395, 6).Value  = 100112001
$ws.Cells.Item(395, 7).Value  = "Berenjena"
$ws.Cells.Item(395, 8).Value  = "Sin especificar"
$ws.Cells.Item(395, 9).Value  = "Primera"
$ws.Cells.Item(395, 10).Value = 124
$ws.Cells.Item(395, 11).Value = 11000
$ws.Cells.Item(395, 12).Value = 12000
$ws.Cells.Item(395, 13).Value = 11500
$ws.Cells.Item(395, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(395, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(395, 16).Value = 230
$ws.Cells.Item(395, 17).Value = 50
$ws.Cells.Item(395, 18).Value = "Hortaliza"
